$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3096085409252669
$wsSummary.Range("C2").Value = 0.06097560975609756
$wsSummary.Range("D2").Value = 0.8928571428571429
$wsSummary.Range("E2").Value = 0.1141552511415525
$wsSummary.Range("F2").Value = 0.239463601532567
$wsSummary.Range("G2").Value = 0.5855855855855856
$wsSummary.Range("H2").Value = 0.7602327447833066
$wsSummary.Range("I2").Value = 25
$wsSummary.Range("J2").Value = 385
$wsSummary.Range("K2").Value = 149
$wsSummary.Range("L2").Value = 3

# ---------------------------------------------------------------
# Sheet: Classification Report
# ---------------------------------------------------------------
$wsClassification = $wb.Worksheets.Item("Classification Report")

$wsClassification.Range("B2").Value = 0.9802631578947368
$wsClassification.Range("C2").Value = 0.2790262172284644
$wsClassification.Range("D2").Value = 0.434402332361516

$wsClassification.Range("B3").Value = 0.06097560975609756
$wsClassification.Range("C3").Value = 0.8928571428571429
$wsClassification.Range("D3").Value = 0.1141552511415525

$wsClassification.Range("B4").Value = 0.3096085409252669
$wsClassification.Range("C4").Value = 0.3096085409252669
$wsClassification.Range("D4").Value = 0.3096085409252669
$wsClassification.Range("E4").Value = 0.3096085409252669

$wsClassification.Range("B5").Value = 0.5206193838254172
$wsClassification.Range("C5").Value = 0.5859416800428037
$wsClassification.Range("D5").Value = 0.2742787917515342

$wsClassification.Range("B6").Value = 0.9344623547846268
$wsClassification.Range("C6").Value = 0.3096085409252669
$wsClassification.Range("D6").Value = 0.4184469617669271

# ---------------------------------------------------------------
# Sheet: Confusion Matrix
# ---------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

$wsConfusion.Range("B2").Value = 149
$wsConfusion.Range("C2").Value = 385

$wsConfusion.Range("B3").Value = 3
$wsConfusion.Range("C3").Value = 25
